$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 601.7
$ws.Range("I6").Value = 557.3333
$ws.Range("K6").Value = 1671.9999
$ws.Range("M6").Value = -1559.9999

$ws.Range("H8").Value = 495
$ws.Range("I8").Value = 495
$ws.Range("K8").Value = 1485
$ws.Range("M8").Value = -1346

$ws.Range("H10").Value = 200
$ws.Range("J10").Value = 200
$ws.Range("L10").Value = 200
$ws.Range("N10").Value = -786

$ws.Range("H12").Value = 4696.391
$ws.Range("I12").Value = 5105.095
$ws.Range("J12").Value = 405
$ws.Range("K12").Value = 5105.095
$ws.Range("L12").Value = 405
$ws.Range("M12").Value = -4935.095
$ws.Range("N12").Value = -745

$ws.Range("H15").Value = 624.2353000000001
$ws.Range("I15").Value = 624.2353000000001
$ws.Range("K15").Value = 1872.7059
$ws.Range("M15").Value = -1703.7059

$ws.Range("H17").Value = 126656.125
$ws.Range("J17").Value = 126656.125
$ws.Range("L17").Value = 379968.375
$ws.Range("N17").Value = -380304.375

$ws.Range("H19").Value = 1837.1
$ws.Range("I19").Value = 1823.875
$ws.Range("J19").Value = 1890
$ws.Range("K19").Value = 1823.875
$ws.Range("L19").Value = 1890
$ws.Range("M19").Value = -1648.875
$ws.Range("N19").Value = -2240

$ws.Range("H40").Value = 3177168.2
$ws.Range("I40").Value = 2105.5
$ws.Range("J40").Value = 8550352
$ws.Range("K40").Value = 2105.5
$ws.Range("L40").Value = 8550352
$ws.Range("M40").Value = -1930.5
$ws.Range("N40").Value = -8550702

$ws.Range("H106").Value = 2930.6155
$ws.Range("I106").Value = 2807.125
$ws.Range("K106").Value = 2807.125
$ws.Range("M106").Value = -2176.125

$ws.Range("H116").Value = 9284.559999999999
$ws.Range("I116").Value = 10339.333
$ws.Range("K116").Value = 10339.333
$ws.Range("M116").Value = -6897.333000000001

$ws.Range("H127").Value = 1669
$ws.Range("I127").Value = 1464.5714
$ws.Range("K127").Value = 4393.7142
$ws.Range("M127").Value = 566.2857999999997

$ws.Range("H138").Value = 3321.1936
$ws.Range("I138").Value = 2857
$ws.Range("K138").Value = 8571
$ws.Range("M138").Value = -3431

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 53300
$ws.Range("I110").Value = 61912.41
$ws.Range("K110").Value = 61912.41
$ws.Range("M110").Value = -59867.41

$ws.Range("H122").Value = 4500.6875
$ws.Range("I122").Value = 2334.75
$ws.Range("K122").Value = 7004.25
$ws.Range("M122").Value = -4554.25

$ws.Range("H124").Value = 51236
$ws.Range("J124").Value = 51236
$ws.Range("L124").Value = 51236
$ws.Range("N124").Value = -61056

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 43328.668
$ws.Range("J81").Value = 43328.668
$ws.Range("L81").Value = 43328.668
$ws.Range("N81").Value = -45450.668

$ws.Range("H84").Value = 43328.668
$ws.Range("J84").Value = 43328.668
$ws.Range("L84").Value = 129986.004
$ws.Range("N84").Value = -140594.004

$ws.Range("H99").Value = 1749.5
$ws.Range("I99").Value = 1749.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1749.5
$ws.Range("L99").Value = 0
$ws.Range("N99").Value = -251.5

$ws.Range("H135").Value = 79999.71000000001
$ws.Range("J135").Value = 79999.71000000001
$ws.Range("L135").Value = 79999.71000000001
$ws.Range("N135").Value = -90139.71000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12175.409
$ws.Range("J31").Value = 15633.333
$ws.Range("L31").Value = 15633.333
$ws.Range("N31").Value = -16223.333

$ws.Range("H34").Value = 12175.409
$ws.Range("J34").Value = 15633.333
$ws.Range("L34").Value = 15633.333
$ws.Range("N34").Value = -16037.333

$ws.Range("H105").Value = 2860209.8
$ws.Range("I105").Value = 4001493.8
$ws.Range("K105").Value = 4001493.8
$ws.Range("M105").Value = -3999746.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 4332.6665
$ws.Range("J45").Value = 4332.6665
$ws.Range("L45").Value = 12997.9995
$ws.Range("N45").Value = -14061.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 59949.5
$ws.Range("J15").Value = 59949.5
$ws.Range("L15").Value = 59949.5
$ws.Range("N15").Value = -60525.5

$ws.Range("H70").Value = 5000
$ws.Range("I70").Value = 5000
$ws.Range("K70").Value = 5000
$ws.Range("M70").Value = -4730

$ws.Range("H73").Value = 5000
$ws.Range("I73").Value = 5000
$ws.Range("K73").Value = 5000
$ws.Range("M73").Value = -4064

$ws.Range("H80").Value = 3000
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0

$ws.Range("H81").Value = 59949.5
$ws.Range("J81").Value = 59949.5
$ws.Range("L81").Value = 59949.5
$ws.Range("N81").Value = -61945.5

$ws.Range("H82").Value = 59999
$ws.Range("J82").Value = 59999
$ws.Range("L82").Value = 59999
$ws.Range("N82").Value = -60765

$ws.Range("H83").Value = 3000
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0

$ws.Range("H84").Value = 59949.5
$ws.Range("J84").Value = 59949.5
$ws.Range("L84").Value = 179848.5
$ws.Range("N84").Value = -189832.5

$ws.Range("H85").Value = 59999
$ws.Range("J85").Value = 59999
$ws.Range("L85").Value = 59999
$ws.Range("N85").Value = -62651

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("N92").Value = 0

$ws.Range("H97").Value = 1327.5834
$ws.Range("I97").Value = 661.5714
$ws.Range("J97").Value = 2260
$ws.Range("K97").Value = 661.5714
$ws.Range("L97").Value = 2260
$ws.Range("M97").Value = -165.5714
$ws.Range("N97").Value = -3252

$ws.Range("H132").Value = 7816120.5
$ws.Range("I132").Value = 7816120.5
$ws.Range("K132").Value = 23448361.5
$ws.Range("M132").Value = -23445831.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2641.3333
$ws.Range("I16").Value = 956.5714
$ws.Range("K16").Value = 956.5714
$ws.Range("M16").Value = -786.5714

$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("N36").Value = 0

$ws.Range("H40").Value = 3191.375
$ws.Range("I40").Value = 3233
$ws.Range("K40").Value = 3233
$ws.Range("M40").Value = -3097

$ws.Range("H46").Value = 3449.6
$ws.Range("J46").Value = 4285.143
$ws.Range("L46").Value = 4285.143
$ws.Range("N46").Value = -4661.143

$ws.Range("H100").Value = 8319120.5
$ws.Range("I100").Value = 12476243
$ws.Range("K100").Value = 12476243
$ws.Range("M100").Value = -12475702

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("N115").Value = 0

$ws.Range("H122").Value = 6292.6665
$ws.Range("I122").Value = 6457.2354
$ws.Range("J122").Value = 3495
$ws.Range("K122").Value = 19371.7062
$ws.Range("L122").Value = 10485
$ws.Range("M122").Value = -16921.7062
$ws.Range("N122").Value = -15385

$ws.Range("H141").Value = 99995
$ws.Range("J141").Value = 99995
$ws.Range("L141").Value = 99995
$ws.Range("N141").Value = -110355

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 49666
$ws.Range("J41").Value = 49666
$ws.Range("L41").Value = 49666
$ws.Range("N41").Value = -50446

$ws.Range("H45").Value = 20000
$ws.Range("J45").Value = 20000
$ws.Range("L45").Value = 20000
$ws.Range("N45").Value = -20982

$ws.Range("H48").Value = 100000
$ws.Range("J48").Value = 100000
$ws.Range("L48").Value = 100000
$ws.Range("N48").Value = -101138

$ws.Range("H96").Value = 3029.4285
$ws.Range("I96").Value = 3102.3333
$ws.Range("J96").Value = 2974.75
$ws.Range("K96").Value = 3102.3333
$ws.Range("L96").Value = 2974.75
$ws.Range("M96").Value = -1729.3333
$ws.Range("N96").Value = -5720.75

$ws.Range("H122").Value = 4468.3687
$ws.Range("I122").Value = 5394.7144
$ws.Range("J122").Value = 1874.6
$ws.Range("K122").Value = 16184.1432
$ws.Range("L122").Value = 5623.799999999999
$ws.Range("M122").Value = -13734.1432
$ws.Range("N122").Value = -10523.8

$ws.Range("H136").Value = 14709803
$ws.Range("I136").Value = 16670599
$ws.Range("K136").Value = 50011797
$ws.Range("M136").Value = -50009247

$ws.Range("H141").Value = 76997.5
$ws.Range("I141").Value = 77000
$ws.Range("J141").Value = 76995
$ws.Range("K141").Value = 77000
$ws.Range("L141").Value = 76995
$ws.Range("M141").Value = -71820
$ws.Range("N141").Value = -87355
